$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.201.08"
$ws.Range("E2").Value = "  -0.18%  "

$ws.Range("D3").Value = "3.013.81"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.21"
$ws.Range("E5").Value = "  +1.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.99"
$ws.Range("E6").Value = "  +1.10%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.012.28"
$ws.Range("E8").Value = "  +0.21%  "

$ws.Range("E9").Value = "  -1.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.37"
$ws.Range("E10").Value = "  +10.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("E11").Value = "  +1.23%  "

$ws.Range("E12").Value = "  -0.92%  "

$ws.Range("E13").Value = "  +1.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.39"
$ws.Range("E14").Value = "  -0.36%  "

$ws.Range("E15").Value = "  +2.54%  "

$ws.Range("D16").Value = "3.516.45"
$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.99"
$ws.Range("E17").Value = "  -0.93%  "

$ws.Range("D18").Value = "62.223.13"
$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("D19").Value = "3.012.85"
$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "447.08"
$ws.Range("E20").Value = "  -2.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.20"
$ws.Range("E21").Value = "  +1.59%  "

$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.31"
$ws.Range("E24").Value = "  +0.84%  "

$ws.Range("E25").Value = "  +10.61%  "

$ws.Range("E26").Value = "  +1.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.03"
$ws.Range("E27").Value = "  -1.39%  "

$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("E29").Value = "  +2.32%  "

$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.17"
$ws.Range("E31").Value = "  +2.99%  "

$ws.Range("E32").Value = "  +1.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.46"
$ws.Range("E33").Value = "  -0.81%  "

$ws.Range("E34").Value = "  +0.51%  "

$ws.Range("D35").Value = "0.0₃0849"
$ws.Range("E35").Value = "  +5.15%  "

$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.83"
$ws.Range("E37").Value = "  +1.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.15"
$ws.Range("E38").Value = "  -0.60%  "

$ws.Range("E39").Value = "  -2.34%  "

$ws.Range("E40").Value = "  -1.13%  "

$ws.Range("E41").Value = "  +2.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.123"
$ws.Range("E42").Value = "  +0.32%  "

$ws.Range("E43").Value = "  +6.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.90"
$ws.Range("E44").Value = "  +9.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "393.98"
$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("E46").Value = "  -1.53%  "

$ws.Range("D47").Value = "2.727.38"
$ws.Range("E47").Value = "  -0.25%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.75"
$ws.Range("E48").Value = "  +4.37%  "

$ws.Range("E50").Value = "  -0.46%  "

$ws.Range("E51").Value = "  -1.32%  "
